{"js": "// Apply the \"Added many more features\" edits to the document body.\n// Each entry is an exact, case-sensitive search string paired with its\n// replacement text. Plain-text search (no wildcards) is used so that\n// characters like \"|\" and \"&\" in the source strings are treated literally.\nconst replacements = [\n  {\n    find: \"Play Magic of the Ring Deluxe Free | Review & Rating\",\n    replace: \"Play Magic of the Ring Deluxe Free - Review & Features\",\n  },\n  {\n    find: \"Fantastic graphics that perfectly reflect the game's magical theme\",\n    replace: \"Fantastic graphics that perfectly reflect the magical theme\",\n  },\n  {\n    find: \"Fully responsive on any operating system making it highly compatible for all players\",\n    replace: \"Smooth compatibility on all devices and operating systems\",\n  },\n  {\n    find: \"Free spins and bonus symbols bring more winning opportunities\",\n    replace: \"Free spins with bonus symbols that expand and offer more winning opportunities\",\n  },\n  {\n    find: \"Limited number of free spins offered\",\n    replace: \"Limited number of bonus features\",\n  },\n  {\n    find: \"Gameplay can be too simplistic for some players\",\n    replace: \"Lack of progressive jackpot\",\n  },\n  {\n    find:\n      \"Read our comprehensive review of Magic of the Ring Deluxe online slot game. Play for free and learn all about the game's features and bonuses.\",\n    replace:\n      \"Discover the magic in Magic of the Ring Deluxe. Play for free and enjoy exciting features.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"text\");\n  await context.sync();\n\n  // \"Play Magic of the Ring Deluxe Free | Review & Rating\" occurs twice in\n  // the document (the page heading and the bolded text near the end) and\n  // both instances map to the same replacement text, so replacing every\n  // match found is correct for all entries in this table.\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the document.\n# Each pair is an exact, case-sensitive Find/Replace run over the whole\n# document (wdReplaceAll), so every matching instance is updated \u2014\n# including the title, which appears both as the page heading and again\n# as bolded text near the end with identical old/new wording.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ find = \"Play Magic of the Ring Deluxe Free | Review & Rating\"; repl = \"Play Magic of the Ring Deluxe Free - Review & Features\" },\n  @{ find = \"Fantastic graphics that perfectly reflect the game's magical theme\"; repl = \"Fantastic graphics that perfectly reflect the magical theme\" },\n  @{ find = \"Fully responsive on any operating system making it highly compatible for all players\"; repl = \"Smooth compatibility on all devices and operating systems\" },\n  @{ find = \"Free spins and bonus symbols bring more winning opportunities\"; repl = \"Free spins with bonus symbols that expand and offer more winning opportunities\" },\n  @{ find = \"Limited number of free spins offered\"; repl = \"Limited number of bonus features\" },\n  @{ find = \"Gameplay can be too simplistic for some players\"; repl = \"Lack of progressive jackpot\" },\n  @{ find = \"Read our comprehensive review of Magic of the Ring Deluxe online slot game. Play for free and learn all about the game's features and bonuses.\"; repl = \"Discover the magic in Magic of the Ring Deluxe. Play for free and enjoy exciting features.\" }\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.Text = $pair.find\n  $find.Replacement.Text = $pair.repl\n  $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
